$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 259
$ws1.Range("F5").Value = 298
$ws1.Range("F7").Value = 370
$ws1.Range("F8").Value = 1954
$ws1.Range("F10").Value = 28
$ws1.Range("F11").Value = 26
$ws1.Range("F12").Value = 1585
$ws1.Range("F13").Value = 1585
$ws1.Range("F14").Value = 1313
$ws1.Range("F18").Value = 13
$ws1.Range("F20").Value = 428
$ws1.Range("F23").Value = 135
$ws1.Range("F24").Value = 6923
$ws1.Range("F25").Value = 7488
$ws1.Range("F26").Value = 31
$ws1.Range("F27").Value = 173
$ws1.Range("F29").Value = 49
$ws1.Range("F36").Value = 1364
$ws1.Range("F40").Value = 666
$ws1.Range("F47").Value = 114
$ws1.Range("F48").Value = 129

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 47

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 248
$ws3.Range("F5").Value = 106

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 259
$ws4.Range("F6").Value = 248
$ws4.Range("F7").Value = 106
$ws4.Range("F9").Value = 298
$ws4.Range("F12").Value = 370
$ws4.Range("F13").Value = 1954
$ws4.Range("F14").Value = 28
$ws4.Range("F15").Value = 26
$ws4.Range("F16").Value = 1585
$ws4.Range("F17").Value = 1585
$ws4.Range("F20").Value = 428
$ws4.Range("F22").Value = 135
$ws4.Range("F23").Value = 47
$ws4.Range("F24").Value = 6923
$ws4.Range("F25").Value = 7488
$ws4.Range("F26").Value = 31
$ws4.Range("F31").Value = 1364
$ws4.Range("F37").Value = 666
$ws4.Range("F43").Value = 303
$ws4.Range("F47").Value = 114
